$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 is the b13c6bb2 file
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-13 10:16:34"

# zh-cn sheet: row 3 is the b13c6bb2 file
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-13 10:16:26"

# de-de sheet: row 3 is the b13c6bb2 file
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-13 10:16:34"

# Adjust column widths (Excel auto-resized these columns after the longer
# "Ready for handoff" status text was written into them). The COM width
# model here quantizes ColumnWidth to 1/6ths, so 16.33 is the input that
# lands on the stored width closest to the target (~17.22 chars).
$overview.Range("E:E").ColumnWidth = 16.33
$overview.Range("F:F").ColumnWidth = 16.33
$zhcn.Range("C:C").ColumnWidth = 16.33
$dede.Range("C:C").ColumnWidth = 16.33
